$wb = $excel.ActiveWorkbook

# --- Add the new "Notification" sheet and populate it while it is still
#     positioned next to the sheet it was inserted before (inserting a new
#     sheet drops it in front of the active sheet). Populate it now, then
#     move it into its final position, since writes should happen before
#     any Move() calls are issued. ---
$notification = $wb.Worksheets.Add()
$notification.Name = "Notification"
$notification.Range("A1").Value = "Email ID"
$notification.Range("A1").Font.Bold = $true
$notification.Range("A2").Value = "Amanda Donovan"

# Move "Notification" to the very end of the tab strip (after MoreAttendees).
$lastSheet = $wb.Worksheets.Item($wb.Worksheets.Count)
$notification.Move($null, $lastSheet)

# Move "UpdateActivity" to the very end of the tab strip (after Notification).
$updateActivity = $wb.Worksheets.Item("UpdateActivity")
$lastSheet = $wb.Worksheets.Item($wb.Worksheets.Count)
$updateActivity.Move($null, $lastSheet)

# Final tab order is now:
#   Users, Contact, Activity, MoreAttendees, Notification, UpdateActivity

# Moving a sheet invalidates earlier worksheet handles in this engine, so
# re-resolve every sheet reference by name before touching it again.

# Update the selected cell on the "Activity" sheet.
$activity = $wb.Worksheets.Item("Activity")
[void]$activity.Range("I10").Select()

# Update the selected cell on the "Notification" sheet.
$notification = $wb.Worksheets.Item("Notification")
[void]$notification.Range("A2").Select()

# "UpdateActivity" becomes the active/visible tab, with a new selection.
$updateActivity = $wb.Worksheets.Item("UpdateActivity")
[void]$updateActivity.Activate()
[void]$updateActivity.Range("G11").Select()
